# bioretention-cells: fix parameterisation bug
# "conductivity_mm.per.hour" row values on the bioretention_cell sheet did not
# line up with their column's scenario name (e.g. the "3.6mm.per.hour_..."
# columns contained 3600 instead of 3.6). Fix those six cells.

$wb = $excel.ActiveWorkbook

# --- 1. bioretention_cell sheet: fix conductivity_mm.per.hour row (row 11) ---
$bio = $wb.Worksheets.Item("bioretention_cell")

$bio.Cells.Item(11, 5).Value  = 3.6   # E11  (3.6mm.per.hour_mulde_rigole_no-drainage)
$bio.Cells.Item(11, 6).Value  = 3.6   # F11  (3.6mm.per.hour_mulde_rigole_with-drainage)
$bio.Cells.Item(11, 8).Value  = 36    # H11  (36mm.per.hour_mulde_rigole_no-drainage)
$bio.Cells.Item(11, 9).Value  = 36    # I11  (36mm.per.hour_mulde_rigole_with-drainage)
$bio.Cells.Item(11, 11).Value = 360   # K11  (360mm.per.hour_mulde_rigole_no-drainage)
$bio.Cells.Item(11, 12).Value = 360   # L11  (360mm.per.hour_mulde_rigole_with-drainage)

# --- 2. green_roof sheet: turn the plain-text "reference" URLs (col M) into
#        real hyperlinks, keeping each row's own target address ---
$roof = $wb.Worksheets.Item("green_roof")

$refUrls = @{
  3  = "https://cloud.kompetenz-wasser.de/index.php/f/180243"
  4  = "https://cloud.kompetenz-wasser.de/index.php/f/180244"
  5  = "https://cloud.kompetenz-wasser.de/index.php/f/180245"
  6  = "https://cloud.kompetenz-wasser.de/index.php/f/180246"
  7  = "https://cloud.kompetenz-wasser.de/index.php/f/180247"
  8  = "https://cloud.kompetenz-wasser.de/index.php/f/180248"
  9  = "https://cloud.kompetenz-wasser.de/index.php/f/180249"
  10 = "https://cloud.kompetenz-wasser.de/index.php/f/180250"
  11 = "https://cloud.kompetenz-wasser.de/index.php/f/180251"
  12 = "https://cloud.kompetenz-wasser.de/index.php/f/180252"
  13 = "https://cloud.kompetenz-wasser.de/index.php/f/180253"
  14 = "https://cloud.kompetenz-wasser.de/index.php/f/180254"
  15 = "https://cloud.kompetenz-wasser.de/index.php/f/180255"
}

$displayText = "https://cloud.kompetenz-wasser.de/index.php/f/180243"

$rowOrder = @(3, 5, 7, 9, 11, 13, 15, 4, 6, 8, 10, 12, 14)
foreach ($r in $rowOrder) {
  $roof.Hyperlinks.Add($roof.Cells.Item($r, 13), $refUrls[$r], "", "", $displayText) | Out-Null
  $roof.Cells.Item($r, 13).Style = "Link"
}

$wb.Save()
